$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inflationbreakdown")

# Append a new data row (row 19) below the existing data (which ends at row 18),
# continuing the monthly date series: 2025-06-01 (serial 45809).
# Copy the formatting (date number format + font) from the last existing date
# cell (A18) so the new date cell matches the style of the rest of column A.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A19").Value = 45809
$ws.Range("B19").Value = 0.26334
$ws.Range("C19").Value = 0.19573
$ws.Range("D19").Value = 0.4332
$ws.Range("E19").Value = 0.10773
$ws.Range("F19").Value = 0.0479
